# "Consolidated Data" sheet: move the overflow metrics that were stuck
# out in columns E:M of rows 2-4 down into their own rows (20-22),
# re-based to start at column A, for complete consolidation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consolidated Data")

# Column A of the new rows holds a "Rat" value that can look numeric
# (e.g. "4.6"); force it to stay text like it was in its original cell.
$ratCol = $ws.Range("A20:A22")
$ratCol.NumberFormat = "@"

# Row 2 (E2:M2) -> Row 20 (A20:I20)
$ws.Range("A20").Value = "4.6"
$ws.Range("B20").Value = 4166
$ws.Range("C20").Value = 4488.02
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 29.5
$ws.Range("G20").Value = 31.78
$ws.Range("H20").Value = 4195.5
$ws.Range("I20").Value = 4519.8

# Row 3 (E3:M3) -> Row 21 (A21:I21)
$ws.Range("A21").Value = "4.9"
$ws.Range("B21").Value = 6266.85
$ws.Range("C21").Value = 6750.95
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 20
$ws.Range("G21").Value = 21.54
$ws.Range("H21").Value = 6286.85
$ws.Range("I21").Value = 6772.49

# Row 4 (E4:M4) -> Row 22 (A22:I22)
$ws.Range("A22").Value = "NA"
$ws.Range("B22").Value = 596
$ws.Range("C22").Value = 641.8
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 30
$ws.Range("G22").Value = 32.31
$ws.Range("H22").Value = 626
$ws.Range("I22").Value = 674.1099999999999

# Drop the formatting we used to keep the "Rat" values textual, so the
# new cells end up plain/unstyled like the rest of the sheet.
$ratCol.ClearFormats()

# The values now live in A20:I22, so clear the old E:M cells on rows 2-4.
$ws.Range("E2:M2").Clear()
$ws.Range("E3:M3").Clear()
$ws.Range("E4:M4").Clear()
